$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("min hjemmeside")
$tc = $r.Font.TextColor
$tc.ObjectThemeColor = 4
$tc.TintAndShade = 0.6
